$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 is the "Aggregate on month" task. The Developer (E) and Status (F)
# columns were left blank - fill them in.
$ws.Range("E15").Value = "Mohandas"
$ws.Range("F15").Value = "[20-May]Still working"

# F15 is a brand new cell in this row; give it the same formatting already
# used by the rest of the table (left/top aligned, wrapped text, no border)
# by copying it from its neighbour E15 rather than rebuilding it property
# by property.
$ws.Range("E15").Copy()
$ws.Range("F15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Leave the selection on the cell that was just edited.
$ws.Range("E15").Select()
